# Lecture partielle de l'EDT M1 MIAGE.
# Update the schedule dates (shifted forward ~3 years) and the corresponding
# day-of-week labels on the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Thursday 2023-03-09 -> Monday 2026-03-09
$ws.Range("A2").Value = 46090
$ws.Range("B2").Value = "lundi"

# Row 4: Saturday 2023-03-11 -> Wednesday 2026-03-11
$ws.Range("A4").Value = 46092
$ws.Range("B4").Value = "mercredi"

# Row 7: Thursday 2023-03-16 -> Monday 2026-03-16
$ws.Range("A7").Value = 46097
$ws.Range("B7").Value = "lundi"

# Row 10: Thursday 2023-05-11 -> Monday 2026-05-11
$ws.Range("A10").Value = 46153
$ws.Range("B10").Value = "lundi"
